$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column B width (closest achievable to target 206.85546875)
$ws.Columns.Item(2).ColumnWidth = 206

# Row 174
$ws.Cells.Item(174,1).Value = "～うちに"
$ws.Cells.Item(174,2).Value = "Trong lúc`n* Thực hiện một hành động trước khi một trạng thái khác thay đổi`n- 日本にいるうちに一度富士山に登ってみたい。`nTrong lúc còn đang ở Nhật, tôi muốn thử leo núi Phú Sĩ một lần.`n* Trong khoảng thời gian một hành động, sự việc xảy ra, có một trạng thái thay đổi,`nkhông mang ý chí`n- 音楽を聞いているうちに眠ってしまった。`n Khi đang nghe nhạc, tôi ngủ thiếp đi mất."
$ws.Cells.Item(174,1).Interior.Color = 65535
$ws.Cells.Item(174,2).WrapText = $true
$ws.Rows.Item(174).RowHeight = 120

# Row 175
$ws.Cells.Item(175,1).Value = "間 (あいだ) … / ～間に(あいだに)… "
$ws.Cells.Item(175,2).Value = "Trong khoảng thời gian`n* ～間… : Khi một hành động xảy ra trong một thời gian dài, một hành động khác`ncũng xảy ra kéo dài đồng thời trong khoảng thời gian đó`n- お母さんが昼寝をしている間、子どもたちはテレビを見ていた。`n Trong khi người mẹ đang ngủ trưa thì bọn trẻ ngồi xem tivi.`n* ～間に: Khi một hành động xảy ra trong một thời gian dài, một hành động khác`nxảy ra chỉ tại một thời điểm trong khoảng thời gian đó`n- 私が旅行で留守の間に、庭に草がたくさん生えてしまった。`n Trong lúc tôi vắng nhà đi du lịch, cỏ dại đã mọc đầy vườn."
$ws.Cells.Item(175,2).WrapText = $true
$ws.Rows.Item(175).RowHeight = 135

# Row 176
$ws.Cells.Item(176,1).Value = "～てからでないと… / ～てからでなければ… "
$ws.Cells.Item(176,2).Value = "Nếu không phải là sau khi ～ thì không thể…`n* Cho đến khi việc gì đó được thực hiện xong, thì một việc khác cũng vẫn chưa làm`nđược. Câu được dùng ở dạng phủ định`n- 運転免許を取ってからでなければ車を運転してはいけない。`n Nếu chưa có bằng thì không được phép lái xe.`n- 病気が治ってからでなければ激しい運動は無理だ。`n Nếu chưa khỏi ốm thì chưa thể vận động mạnh được"
$ws.Cells.Item(176,2).WrapText = $true
$ws.Rows.Item(176).RowHeight = 105

# Row 177
$ws.Cells.Item(177,1).Value = "～ところだ / ～ところ (+trợ từ)"
$ws.Cells.Item(177,2).Value = "Vào thời điểm`n* Khoảnh khắc ngay trước, khoảnh khắc ngay sau, hay vào chính giữa khoảnh khắc`nđó`n- ロケットは間もなく飛び立つところです。緊張の瞬間です。`n Quả tên lửa chuẩn bị được phóng đi rồi. Khoảnh khắc này thật là hồi hộp.`n- 試験中、となりの人の答えを見ているところを先生に注意された。`n Trong giờ thi, tôi bị cô giáo nhắc đúng lúc đang nhìn bài bạn bên cạnh"
$ws.Cells.Item(177,2).WrapText = $true
$ws.Rows.Item(177).RowHeight = 105

# Row 178
$ws.Cells.Item(178,1).Value = "N の/ Vる /Vた/・・・＋ ～とおりだ / ～とろり（に）`nN ＋ ～どおりだ /～どおり（に）"
$ws.Cells.Item(178,2).Value = "Dựa theo`n* Thực hiện một hành động giống như hướng dẫn có sẵn, hoặc giống như cách`nngười khác làm trước`n- 交番で 教えてもらったとおりに歩いていったので、迷わず会場に着いた。`n Nhờ đi theo chỉ dẫn của cảnh sát, tôi đã đến được hội trường mà không bị lạc.`n- サッカーの試合の結果はわたしたちの期待どおりだった。`n Kết quả trận bóng đúng như những gì tôi mong đợi"
$ws.Cells.Item(178,1).WrapText = $true
$ws.Cells.Item(178,2).WrapText = $true
$ws.Rows.Item(178).RowHeight = 105

# Row 179
$ws.Cells.Item(179,1).Value = "～によって… ～によっては…"
$ws.Cells.Item(179,2).Value = "Tùy vào, tùy theo`n* Thể hiện trạng thái không nhất định mà thay đổi tùy thuộc vào một yếu tố nào`nđó`n- わたしの帰宅時間は毎日違う。日によって 夜中になることもある。`n Giờ giấc về nhà của tôi mỗi ngày một khác nhau. Tùy hôm, có lúc đêm muộn mới về được.`n- 国によって習慣が違う。`n Tùy từng nước, tập quán lại khác nhau."
$ws.Cells.Item(179,2).WrapText = $true
$ws.Rows.Item(179).RowHeight = 105

# Row 180
$ws.Cells.Item(180,1).Value = "N の/ Vる ＋ ～たびに…"
$ws.Cells.Item(180,2).Value = "Mỗi lần`n* Mỗi khi điều gì xảy ra thì có điều khác cũng xảy ra theo. Nhấn mạnh sự lặp lại của`nhành động, nhưng không dùng cho những hiện tượng tự nhiên diễn ra hàng ngày`n- この地方は台風が来るたびに大水の害が起こる。`n Khu vực này mỗi khi có bão lại phải chịu thiệt hại do ngập lụt.`n- このチームは試合のたびに強くなっていく。`n Đội (bóng) nay mạnh dần lên sau mỗi trận đấu."
$ws.Cells.Item(180,2).WrapText = $true
$ws.Rows.Item(180).RowHeight = 105

# Row 181
$ws.Cells.Item(181,1).Value = "（～ば）～ほど /（～なら）～ほど/ ～ほど"
$ws.Cells.Item(181,2).Value = "Càng - càng`n* Phía sau ほど ・・・ là một vế thể hiện sự thay đổi tương ứng, tỉ lệ với vế trước`n- ものが増えれば増えるほど整理が大変になる。`n Càng nhiều đồ thì việc sắp xếp càng khó hơn.`n- 忙しい人ほど時間の使い方が上手だ。`n Những người càng bận rộn thì càng giỏi quản lý thời gian hơn."
$ws.Cells.Item(181,2).WrapText = $true
$ws.Rows.Item(181).RowHeight = 90

# Row 182
$ws.Cells.Item(182,1).Value = "Nの/Vる/Vた + ～ついでに"
$ws.Cells.Item(182,2).Value = "Tiện thể`n* Nhân lúc thực hiện một hành động này thì thì thực hiện thêm một hành động`nkhác`n- 散歩のついでにこのはがきをポストに出してきた。`n Nhân tiện ra ngoài đi dạo, tôi đã đi gửi bưu thiếp.`n- インタネットで本を注文したついでに新しく出たDVDも調べた。`n Nhân tiện đặt mua sách trên mạng, tôi cũng thử tìm hiểu đĩa DVD mới ra."
$ws.Cells.Item(182,2).WrapText = $true
$ws.Rows.Item(182).RowHeight = 105

# Row 183
$ws.Cells.Item(183,1).Value = "N / Vる・ない (Naだ -な) ＋ ～くらい・～ぐらい・～ほど… "
$ws.Cells.Item(183,2).Value = "Tới mức`n* Biểu thị mức độ của một trạng thái sự vật, hiện tượng bằng cách đưa ra một sự vật`nhiện tượng khác`n- 天気予報によると、今日は台風ぐらい風が吹くそうだ。`n Theo dự báo thời tiết, hôm nay có gió mạnh như bão.`n- さっき地震があった。本棚が倒れるかと思うほど激しく揺れた。`n Vừa mới có động đất. Mặt đất rung mạnh đến mức tôi tưởng kệ sách bị đổ."
$ws.Cells.Item(183,2).WrapText = $true
$ws.Rows.Item(183).RowHeight = 105

# Final selection / view state (best effort)
$ws.Range("B163").Select()
